$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 92
$ws1.Range("F4").Value = 7253
$ws1.Range("F5").Value = 266
$ws1.Range("F7").Value = 3741
$ws1.Range("F9").Value = 536
$ws1.Range("F12").Value = 99

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 92
$ws4.Range("F5").Value = 7253
$ws4.Range("F7").Value = 266
$ws4.Range("F9").Value = 3741
$ws4.Range("F11").Value = 536
$ws4.Range("F14").Value = 99
